$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), Q (Unidad de comercializacion),
# S (Precio $/Kg). Column T (Kg/unidad) is unchanged (14 everywhere).

$rows = @{
    2 = @{ D = 44162; M = 120; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 }
    3 = @{ D = 44351; M = 300; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos empedrada"; S = 714 }
    4 = @{ D = 44309; M = 300; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 }
    5 = @{ D = 44176; M = 250; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 }
    6 = @{ D = 44400; M = 100; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos";            S = 714 }
    7 = @{ D = 44397; M = 60;  N = 11000; O = 11000; P = 11000; Q = "`$/caja 14 kilos";            S = 786 }
    8 = @{ D = 44208; M = 210; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos empedrada"; S = 714 }
    9 = @{ D = 44491; M = 180; N = 9000;  O = 9000;  P = 9000;  Q = "`$/caja 14 kilos empedrada"; S = 643 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("S$r").Value = $vals.S
}
